$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 quantity cells C2:AG2
$ws.Range("C2").Value = 1.999
$ws.Range("D2:AG2").Value = 1.99

# Update rate cell used by AL formulas
$ws.Range("AP2").Value = 42

# Update the two date cells (AT2, AU2)
$ws.Range("AT2").Value = 43344
$ws.Range("AU2").Value = 43390

# Restore the view to the top-left of the sheet and move the selection
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M7").Select()
